$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerial = 45917
$progress = 0.9052869546929831

$rows = @(
    @{ Row = 52; A = "G2"; B = "Workout" },
    @{ Row = 53; A = "G3"; B = "Eat Healthy" },
    @{ Row = 54; A = "G4"; B = "Read Book" },
    @{ Row = 55; A = "G5"; B = "Investment Plan" },
    @{ Row = 56; A = "G6"; B = "Spend 10 Hours without phone" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $cDate = $ws.Cells.Item($rowNum, 3)
    $cDate.Value = $dateSerial
    $cDate.NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($rowNum, 4).Value = $progress
    $ws.Cells.Item($rowNum, 5).Value = 0
    $ws.Cells.Item($rowNum, 6).Value = -0.01
}
